# fixes per #20; regenerated files
#
# Two textual fixes inside the <head><m>...</m></head> markup-annotated
# transcription text:
#
#   1. "E<corr><del>e</del></corr>au de sel armoniac"
#        -> "Eeau de sel armoniac"
#      i.e. drop the <corr><del> ... </del></corr> correction markup
#      around the redundant "e", leaving the underlying reading "eau"
#      intact as plain text (same run formatting as its neighbours).
#
#   2. "recuit) <del>c</del> quand tu le vouldras destremper d'eau..."
#        -> "recuit) <del>C</del> quand tu le vouldras destremper d'eau..."
#      i.e. capitalize the "c" that starts the new sentence/paragraph.
#
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Fix 1: "<corr><del>e</del></corr>au" -> "eau"
# ---------------------------------------------------------------------
# The three plain-text runs around the markup ("E", "e" and "au") all
# share identical run formatting (color 000000, rtl 0), so a naive
# delete of the markup runs causes the host to auto-merge *every*
# touching run of identical formatting together (not just "e"+"au", but
# also the preceding "E" and the following " de sel armoniac"). To keep
# the run boundaries exactly as in the target (only "e"+"au" combine
# into "eau", while "E" and " de sel armoniac" remain separate runs) we
# temporarily perturb the font color of the runs on *both* sides of the
# edit so they cannot be merged into, perform the markup deletions, then
# restore their original color afterwards.

$text = $d.Content.Text

# Temporarily mark the "E" run so it won't absorb "eau" once the
# <corr><del> run between them disappears.
$idxE = $text.IndexOf("E<corr><del>")
$rE = $d.Range($idxE, $idxE + 1)
$rE.Font.Color = 2

# Temporarily mark the " de sel armoniac" run so it won't absorb "eau"
# once the </del></corr> run between them disappears.
$text = $d.Content.Text
$idxTail = $text.IndexOf(" de sel armoniac")
$lenTail = " de sel armoniac".Length
$rTail = $d.Range($idxTail, $idxTail + $lenTail)
$rTail.Font.Color = 1

# Delete "</del></corr>" (do the rightmost markup span first so the
# left span's offset stays valid).
$text = $d.Content.Text
$idxClose = $text.IndexOf("</del></corr>")
$lenClose = "</del></corr>".Length
$rClose = $d.Range($idxClose, $idxClose + $lenClose)
$rClose.Delete()

# Delete "<corr><del>"
$text = $d.Content.Text
$idxOpen = $text.IndexOf("<corr><del>")
$lenOpen = "<corr><del>".Length
$rOpen = $d.Range($idxOpen, $idxOpen + $lenOpen)
$rOpen.Delete()

# Restore the "E" run's original color (black).
$text = $d.Content.Text
$idxE2 = $text.IndexOf("Eeau")
$rE2 = $d.Range($idxE2, $idxE2 + 1)
$rE2.Font.Color = 0

# Restore the " de sel armoniac" run's original color (black).
$text = $d.Content.Text
$idxTail2 = $text.IndexOf(" de sel armoniac")
$rTail2 = $d.Range($idxTail2, $idxTail2 + $lenTail)
$rTail2.Font.Color = 0

# ---------------------------------------------------------------------
# Fix 2: capitalize "c" -> "C" in "recuit) <del>c</del> quand..."
# ---------------------------------------------------------------------
$text = $d.Content.Text
$idxDel = $text.IndexOf("recuit) <del>c</del>")
$idxC = $idxDel + "recuit) <del>".Length
$rC = $d.Range($idxC, $idxC + 1)
$rC.Text = "C"
